# Update the dSF column (F) values for the affected rows to reflect
# the repulled data / push-all-data / mean calculation changes.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F4").Value = 0
$ws.Range("F7").Value = -6
$ws.Range("F8").Value = 3
$ws.Range("F10").Value = 1
$ws.Range("F11").Value = 1
$ws.Range("F12").Value = 4
